$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that just get a smaller D value (B unchanged, E = D-B shared formula auto-recalcs) ---
$ws.Range("D3").Value2 = 5
$ws.Range("D5").Value2 = 5
$ws.Range("D6").Value2 = 4
$ws.Range("D10").Value2 = 5
$ws.Range("D12").Value2 = 5
$ws.Range("D13").Value2 = 5
$ws.Range("D16").Value2 = 5
$ws.Range("D17").Value2 = 10
$ws.Range("D18").Value2 = 9
$ws.Range("D22").Value2 = 9
$ws.Range("D23").Value2 = 10
$ws.Range("D25").Value2 = 9
$ws.Range("D27").Value2 = 9
$ws.Range("D29").Value2 = 10
$ws.Range("D30").Value2 = 8

# Row 27 also loses its "takes long" note in column F
$ws.Range("F27").ClearContents()

# --- Row 8: was an error row ("err" / #VALUE! / "cherry reduction at supress"), now a plain numeric row ---
$ws.Range("D8").Value2 = 5
$ws.Range("F8").ClearContents()

# --- Row 24: was an error row, now plain numeric ---
$ws.Range("D24").Value2 = 10
$ws.Range("F24").ClearContents()

# --- Row 28: was an error row, now plain numeric ---
$ws.Range("D28").Value2 = 10
$ws.Range("F28").ClearContents()

# --- Rows 14, 15, 19, 21: were plain numeric rows, now become error rows ---
$ws.Range("D14").Value2 = "err"
$ws.Range("F14").Value2 = "singleton suppress"

$ws.Range("D15").Value2 = "err"
$ws.Range("F15").Value2 = "singleton suppress"

$ws.Range("D19").Value2 = "err"
$ws.Range("F19").Value2 = "singleton suppress"

$ws.Range("D21").Value2 = "err"
$ws.Range("F21").Value2 = "singleton suppress"

# --- Selection moved from B8 to H35 ---
$ws.Range("H35").Select()
